$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns at B and C, shifting old B (Weight_before) and
# old C (Weight_after) to D and E respectively.
$ws.Range("B:C").Insert()

# Update header row
$ws.Range("A1").Value = "Sample"
$ws.Range("B1").Value = "Chert_type"
$ws.Range("C1").Value = "Chert_tool"
$ws.Range("D1").Value = "Weight_before_[mg]"
$ws.Range("E1").Value = "Weight_after_[mg]"

# Copy header style (bold/centered) from A1 onto the two new header cells
$ws.Range("A1").Copy()
$ws.Range("B1:C1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Sample IDs, Chert type, Chert tool number
$samples = @("ISEA-EX1","ISEA-EX2","ISEA-EX3","ISEA-EX4","ISEA-EX5","ISEA-EX6",
             "ISEA-EX7","ISEA-EX8","ISEA-EX9","ISEA-EX10","ISEA-EX11","ISEA-EX12")
$chertType = @("A","B","A","B","A","B","A","B","A","B","A","B")
$chertTool = @(1,1,2,2,3,3,4,4,5,5,6,6)

for ($i = 0; $i -lt 12; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $samples[$i]
    $ws.Cells.Item($row, 2).Value = $chertType[$i]
    $ws.Cells.Item($row, 3).Value = $chertTool[$i]
}

$ws.Range("A1:E13").EntireColumn.AutoFit()
